# Generate Report for Handback
# The a4fc9b30-... handback file has now been processed: its status moves
# from "Ready for handoff" / "The version ... is not the latest" to
# "Handed back: in sync with en-US", the Latest Handback DateTime is
# refreshed, and the old error detail is cleared.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-15 16:45:36"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).AutoFit()

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-15 16:45:43"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).AutoFit()
